# Add 2022-Q3 data:
#  - insert a new worksheet "2022-Q3" right after "总计" (so it becomes the
#    2nd sheet, shifting 2022-Q2 .. 2021-Q1 down by one position each);
#  - populate it with the fund-holdings table for the new quarter;
#  - update the "总计" (summary) sheet with a new leading row for 2022-Q3
#    and renumber the trailing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after the first sheet ("总计").
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row (row 1) - bold/centered header style used by every other
# quarterly sheet in this workbook.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"


# Data columns B..G on rows 2-6 hold text values (fund code / name / size /
# position / position% / held-value) exactly like every other quarterly
# sheet - mark them as Text *before* assigning so Excel doesn't silently
# reinterpret the numeric-looking strings as numbers.
$q3.Range("B2:G6").NumberFormat = "@"

# Row 2 - 006102 浙商丰利增强债券
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "006102"
$q3.Range("C2").Value = "浙商丰利增强债券"
$q3.Range("D2").Value = "93.25"
$q3.Range("E2").Value = "44.13"
$q3.Range("F2").Value = "1.94"
$q3.Range("G2").Value = "1.8090"
$q3.Range("H2").Value = 6

# Row 3 - 010381 浙商智选价值混合A
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "010381"
$q3.Range("C3").Value = "浙商智选价值混合A"
$q3.Range("D3").Value = "4.54"
$q3.Range("E3").Value = "91.07"
$q3.Range("F3").Value = "4.39"
$q3.Range("G3").Value = "0.1993"
$q3.Range("H3").Value = 10

# Row 4 - 010382 浙商智选价值混合C
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "010382"
$q3.Range("C4").Value = "浙商智选价值混合C"
$q3.Range("D4").Value = "2.31"
$q3.Range("E4").Value = "91.07"
$q3.Range("F4").Value = "4.39"
$q3.Range("G4").Value = "0.1014"
$q3.Range("H4").Value = 10

# Row 5 - 163110 申万菱信量化小盘股票（LOF）A
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "163110"
$q3.Range("C5").Value = "申万菱信量化小盘股票（LOF）A"
$q3.Range("D5").Value = "5.04"
$q3.Range("E5").Value = "93.06"
$q3.Range("F5").Value = "0.60"
$q3.Range("G5").Value = "0.0302"
$q3.Range("H5").Value = 4

# Row 6 - 013918 申万菱信量化小盘股票（LOF）C (held value is exactly 0 -> numeric 0,
# matching the convention used elsewhere in this workbook for zero holdings)
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "013918"
$q3.Range("C6").Value = "申万菱信量化小盘股票（LOF）C"
$q3.Range("D6").Value = "0.00"
$q3.Range("E6").Value = "93.06"
$q3.Range("F6").Value = "0.60"
$q3.Range("G6").NumberFormat = "General"
$q3.Range("G6").Value = 0
$q3.Range("H6").Value = 4



# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 right
#    after the header, shifting every existing quarter row down by one.
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 5
$zongji.Range("D2").Value = 2.14


# Renumber the index column (A) for the rows that shifted down so it stays
# a contiguous 0-based sequence: 2022-Q2=1, 2022-Q1=2, 2021-Q4=3, 2021-Q3=4,
# 2021-Q2=5, 2021-Q1=6.
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
$zongji.Range("A6").Value = 4
$zongji.Range("A7").Value = 5
$zongji.Range("A8").Value = 6

